# The weekly update adds one new price observation for "Rabanito" at
# Vega Modelo de Temuco. It is inserted as the new row 6, pushing every
# existing data row (old rows 6..137) down by one (new rows 7..138).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 6 (shifts rows 6:137 down to 7:138,
# carrying their formatting/styles with them).
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new weekly observation.
$ws.Range("A6").Value = 10
$ws.Range("B6").Value = "Vega Modelo de Temuco"
$ws.Range("C6").Value = "La Araucanía"
$ws.Range("D6").Value2 = 45237
$ws.Range("E6").Value = 9
$ws.Range("F6").Value = 300000001
$ws.Range("G6").Value = "Rabanito"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 40
$ws.Range("K6").Value = 9000
$ws.Range("L6").Value = 9000
$ws.Range("M6").Value = 9000
$ws.Range("N6").Value = "`$/docena de paquetes"
$ws.Range("O6").Value = "Provincia de Cautín"
$ws.Range("P6").Value = 750
$ws.Range("Q6").Value = 12
$ws.Range("R6").Value = "Hortaliza"
